# Updated cryptos list on Thu Sep 19 23:16:03 UTC 2024 with GitHub Actions
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for each coin row.
# A leading apostrophe is used for the Price column so Excel stores the
# refreshed value as text (preserving formatting such as trailing zeros,
# thousands-style dot separators, and very small decimals) instead of
# re-interpreting it as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.992.23'
$ws.Range("E2").Value = '  +3.39%  '
$ws.Range("D3").Value = '''2.475.24'
$ws.Range("E3").Value = '  +5.37%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''565.78'
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("D6").Value = '''142.45'
$ws.Range("E6").Value = '  +7.90%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '''0.588'
$ws.Range("E8").Value = '  +0.96%  '
$ws.Range("D9").Value = '''2.472.23'
$ws.Range("E9").Value = '  +5.35%  '
$ws.Range("E10").Value = '  +2.90%  '
$ws.Range("E11").Value = '  +1.18%  '
$ws.Range("E12").Value = '  +1.36%  '
$ws.Range("D13").Value = '''0.352'
$ws.Range("E13").Value = '  +3.75%  '
$ws.Range("D14").Value = '''26.61'
$ws.Range("E14").Value = '  +10.47%  '
$ws.Range("D15").Value = '''2.915.99'
$ws.Range("E15").Value = '  +5.47%  '
$ws.Range("D16").Value = '''62.862.32'
$ws.Range("E16").Value = '  +3.38%  '
$ws.Range("D17").Value = '''0.0000142'
$ws.Range("E17").Value = '  +4.67%  '
$ws.Range("D18").Value = '''2.473.12'
$ws.Range("E18").Value = '  +5.20%  '
$ws.Range("D19").Value = '''11.25'
$ws.Range("E19").Value = '  +4.71%  '
$ws.Range("D20").Value = '''340.42'
$ws.Range("E20").Value = '  +7.76%  '
$ws.Range("E21").Value = '  +3.37%  '
$ws.Range("D22").Value = '''6.83'
$ws.Range("E22").Value = '  +2.48%  '
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").Value = '''65.59'
$ws.Range("E24").Value = '  +1.94%  '
$ws.Range("E25").Value = '  +1.23%  '
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("E27").Value = '  +5.49%  '
$ws.Range("D28").Value = '''8.10'
$ws.Range("E28").Value = '  +1.11%  '
$ws.Range("D29").Value = '''1.38'
$ws.Range("E29").Value = '  +8.68%  '
$ws.Range("D30").Value = '''6.85'
$ws.Range("E30").Value = '  +11.60%  '
$ws.Range("D31").Value = '''1.85'
$ws.Range("E31").Value = '  +5.99%  '
$ws.Range("D32").Value = '''0.0₃0801'
$ws.Range("E32").Value = '  +8.27%  '
$ws.Range("D33").Value = '''176.92'
$ws.Range("E33").Value = '  +3.23%  '
$ws.Range("E34").Value = '  +10.16%  '
$ws.Range("D35").Value = '''0.399'
$ws.Range("E35").Value = '  +3.35%  '
$ws.Range("D36").Value = '''18.83'
$ws.Range("E36").Value = '  +3.98%  '
$ws.Range("D37").Value = '''370.58'
$ws.Range("E37").Value = '  +10.90%  '
$ws.Range("E38").Value = '  +5.38%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("E41").Value = '  +9.40%  '
$ws.Range("D43").Value = '''149.69'
$ws.Range("E43").Value = '  +7.05%  '
$ws.Range("E44").Value = '  +4.78%  '
$ws.Range("D45").Value = '''20.54'
$ws.Range("E45").Value = '  +5.64%  '
$ws.Range("E46").Value = '  +4.65%  '
$ws.Range("E47").Value = '  +0.72%  '
$ws.Range("D48").Value = '''0.0516'
$ws.Range("E48").Value = '  +2.83%  '
$ws.Range("E49").Value = '  +4.15%  '
$ws.Range("E50").Value = '  +1.76%  '
$ws.Range("D51").Value = '''17.97'
$ws.Range("E51").Value = '  +4.08%  '
